$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.00", "582.89") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.920.49'
$ws.Range("E2").Value = '  -2.50%  '
$ws.Range("D3").Value = '3.004.23'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '582.89'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").Value = '146.04'
$ws.Range("E6").Value = '  -5.54%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -2.60%  '
$ws.Range("D9").Value = '3.001.76'
$ws.Range("E9").Value = '  -2.29%  '
$ws.Range("E10").Value = '  -5.17%  '
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -2.42%  '
$ws.Range("E13").Value = '  -4.37%  '
$ws.Range("D14").Value = '34.34'
$ws.Range("E14").Value = '  -6.35%  '
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '3.497.85'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '6.98'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '61.924.86'
$ws.Range("E18").Value = '  -2.45%  '
$ws.Range("D19").Value = '3.008.19'
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("D20").Value = '459.82'
$ws.Range("E20").Value = '  -5.09%  '
$ws.Range("D21").Value = '13.87'
$ws.Range("E21").Value = '  -4.25%  '
$ws.Range("D22").Value = '0.679'
$ws.Range("E22").Value = '  -4.12%  '
$ws.Range("D23").Value = '7.42'
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").Value = '81.52'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  -6.63%  '
$ws.Range("D26").Value = '12.18'
$ws.Range("E26").Value = '  -5.04%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '9.88'
$ws.Range("E28").Value = '  -6.60%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").Value = '2.61'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").Value = '7.09'
$ws.Range("E31").Value = '  -6.38%  '
$ws.Range("E32").Value = '  -6.55%  '
$ws.Range("D33").Value = '27.57'
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("E35").Value = '  -4.48%  '
$ws.Range("D36").Value = '0.0₃0782'
$ws.Range("E36").Value = '  -5.01%  '
$ws.Range("E37").Value = '  -5.15%  '
$ws.Range("E38").Value = '  -5.96%  '
$ws.Range("D39").Value = '50.00'
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("D40").Value = '8.98'
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("D41").Value = '2.89'
$ws.Range("E41").Value = '  -11.02%  '
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("E43").Value = '  -6.69%  '
$ws.Range("D44").Value = '0.0351'
$ws.Range("E44").Value = '  -3.11%  '
$ws.Range("D45").Value = '377.89'
$ws.Range("E45").Value = '  -14.80%  '
$ws.Range("D46").Value = '2.742.68'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '37.21'
$ws.Range("E47").Value = '  -6.65%  '
$ws.Range("D48").Value = '127.66'
$ws.Range("E48").Value = '  -3.96%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("E51").Value = '  -3.65%  '

# Restore original (default) style on column D now that values are set,
# so the cells keep the same formatting as before the edit.
$ws.Range("D2:D51").Style = "Normal"
